$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.234
$ws.Range("C7").Value = -13.547
$ws.Range("B9").Value = 5.703
$ws.Range("C12").Value = -11.108
$ws.Range("B13").Value = 5.601
$ws.Range("C14").Value = -12.865
$ws.Range("E15").Value = 16.346
$ws.Range("B16").Value = 5.534
$ws.Range("B18").Value = 5.01
$ws.Range("C19").Value = -12.059
$ws.Range("B20").Value = 7.4
$ws.Range("B26").Value = 6.331999999999999
$ws.Range("C26").Value = -12.714
$ws.Range("B27").Value = 6.628
$ws.Range("C27").Value = -12.723
$ws.Range("E28").Value = 16.872
$ws.Range("B29").Value = 5.143
$ws.Range("C29").Value = -11.056
$ws.Range("E33").Value = 17.536
$ws.Range("B35").Value = 7.537000000000001
$ws.Range("E35").Value = 16.527
$ws.Range("B36").Value = 7.971999999999999
$ws.Range("C37").Value = -13.238
$ws.Range("C38").Value = -13.123
$ws.Range("E38").Value = 15.912
$ws.Range("E43").Value = 17.038
$ws.Range("E44").Value = 16.769
$ws.Range("B45").Value = 5.736
$ws.Range("E45").Value = 16.535
$ws.Range("C47").Value = -12.953
$ws.Range("E47").Value = 15.889
$ws.Range("C51").Value = -11.036
$ws.Range("E51").Value = 17.01
$ws.Range("C52").Value = -11.621
$ws.Range("E54").Value = 16.42
$ws.Range("B55").Value = 5.552
$ws.Range("C55").Value = -13.244
$ws.Range("B57").Value = 5.147
$ws.Range("E57").Value = 16.48
$ws.Range("E62").Value = 16.323
$ws.Range("E63").Value = 17.389
$ws.Range("E67").Value = 17.248
$ws.Range("B69").Value = 5.226000000000001
$ws.Range("C69").Value = -10.732
$ws.Range("C70").Value = -11.959
$ws.Range("E70").Value = 17.482
$ws.Range("B76").Value = 5.891
$ws.Range("C76").Value = -12.587
$ws.Range("B78").Value = 7.523999999999999
$ws.Range("C81").Value = -13.237
$ws.Range("E81").Value = 16.643
$ws.Range("B82").Value = 5.061
$ws.Range("B83").Value = 5
$ws.Range("C83").Value = -13.98
$ws.Range("E88").Value = 16.287
$ws.Range("B93").Value = 5.976000000000001
$ws.Range("C94").Value = -10.869
$ws.Range("E96").Value = 16.196
$ws.Range("B97").Value = 5.052
$ws.Range("E99").Value = 16.244
$ws.Range("C100").Value = -12.951
$ws.Range("C102").Value = -13.306
